$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 3.921224538834813
$ws.Range("D2").Value = 10.33175560687762
$ws.Range("E2").Value = 14.14541887378391
$ws.Range("F2").Value = 27.50182019754018
$ws.Range("G2").Value = 25.37668772468433
$ws.Range("H2").Value = 13.56558100888579
$ws.Range("I2").Value = 18.88446172054069
$ws.Range("J2").Value = 9.724607741583178
$ws.Range("K2").Value = 16.36480086728358
$ws.Range("O2").Value = 20.17521814148836

# Row 3
$ws.Range("C3").Value = 3.769521156566517
$ws.Range("D3").Value = 10.24728290284371
$ws.Range("E3").Value = 14.06860491544868
$ws.Range("F3").Value = 27.62822780942313
$ws.Range("G3").Value = 25.5885374952674
$ws.Range("H3").Value = 13.64075277874275
$ws.Range("I3").Value = 18.99485909703471
$ws.Range("J3").Value = 9.732238681792538
$ws.Range("K3").Value = 15.59808833808076
$ws.Range("O3").Value = 20.31402466094038

# Row 4
$ws.Range("C4").Value = 3.672465037595984
$ws.Range("D4").Value = 10.19646961889086
$ws.Range("E4").Value = 14.02399882712069
$ws.Range("F4").Value = 27.71472700142449
$ws.Range("G4").Value = 25.73069384691065
$ws.Range("H4").Value = 13.68974684905275
$ws.Range("I4").Value = 19.06835994336339
$ws.Range("J4").Value = 9.738524504181113
$ws.Range("K4").Value = 15.10632208220523
$ws.Range("O4").Value = 20.40514592377729

# Row 5
$ws.Range("C5").Value = 3.631966836922608
$ws.Range("D5").Value = 10.17604519987967
$ws.Range("E5").Value = 14.00647937647123
$ws.Range("F5").Value = 27.75219929036929
$ws.Range("G5").Value = 25.79163615331798
$ws.Range("H5").Value = 13.71042609611806
$ws.Range("I5").Value = 19.09974423156537
$ws.Range("J5").Value = 9.741488669719246
$ws.Range("K5").Value = 14.90084397973922
$ws.Range("O5").Value = 20.44375695808005

# Row 6
$ws.Range("C6").Value = 3.625186088731983
$ws.Range("D6").Value = 10.17267130066281
$ws.Range("E6").Value = 14.00361043018565
$ws.Range("F6").Value = 27.75855548538901
$ws.Range("G6").Value = 25.80193678994758
$ws.Range("H6").Value = 13.71390298128276
$ws.Range("I6").Value = 19.10504190810152
$ws.Range("J6").Value = 9.742005189766493
$ws.Range("K6").Value = 14.86642425886728
$ws.Range("O6").Value = 20.45025747874125

# Row 7
$ws.Range("C7").Value = 3.67192264815723
$ws.Range("D7").Value = 10.19619300174851
$ws.Range("E7").Value = 14.02375987121999
$ws.Range("F7").Value = 27.71522337892554
$ws.Range("G7").Value = 25.73150357201687
$ws.Range("H7").Value = 13.69002284651594
$ws.Range("I7").Value = 19.0687774114051
$ws.Range("J7").Value = 9.738562849474816
$ws.Range("K7").Value = 15.10357120825168
$ws.Range("O7").Value = 20.40566066446528

# Row 8
$ws.Range("C8").Value = 3.869748867915773
$ws.Range("D8").Value = 10.30242147084157
$ws.Range("E8").Value = 14.11841158619497
$ws.Range("F8").Value = 27.5435551211712
$ws.Range("G8").Value = 25.44721123539117
$ws.Range("H8").Value = 13.59091117271554
$ws.Range("I8").Value = 18.92133708111048
$ws.Range("J8").Value = 9.726906883074399
$ws.Range("K8").Value = 16.10490843713621
$ws.Range("O8").Value = 20.22185325061334

# Row 9
$ws.Range("C9").Value = 4.225307269550791
$ws.Range("D9").Value = 10.51830241977484
$ws.Range("E9").Value = 14.32363393584636
$ws.Range("F9").Value = 27.27788179471771
$ws.Range("G9").Value = 24.98669011797944
$ws.Range("H9").Value = 13.41907162908514
$ws.Range("I9").Value = 18.677803195653
$ws.Range("J9").Value = 9.716735483961573
$ws.Range("K9").Value = 17.89457083493884
$ws.Range("O9").Value = 19.90833322528667

# Row 10
$ws.Range("C10").Value = 4.465266375492308
$ws.Range("D10").Value = 10.68044327352582
$ws.Range("E10").Value = 14.48539189574181
$ws.Range("F10").Value = 27.12658276840904
$ws.Range("G10").Value = 24.70901319003152
$ws.Range("H10").Value = 13.30654529988185
$ws.Range("I10").Value = 18.5270027807821
$ws.Range("J10").Value = 9.716974513741858
$ws.Range("K10").Value = 19.09544002879056
$ws.Range("O10").Value = 19.70682612443444

# Row 11
$ws.Range("C11").Value = 4.569559659042902
$ws.Range("D11").Value = 10.75473844097327
$ws.Range("E11").Value = 14.56114165841566
$ws.Range("F11").Value = 27.06741035260248
$ws.Range("G11").Value = 24.59619783636018
$ws.Range("H11").Value = 13.25833437486849
$ws.Range("I11").Value = 18.46457718923851
$ws.Range("J11").Value = 9.718750898245553
$ws.Range("K11").Value = 19.61575993144636
$ws.Range("O11").Value = 19.62146880934766

# Row 12
$ws.Range("C12").Value = 4.608336150573902
$ws.Range("D12").Value = 10.78292942104152
$ws.Range("E12").Value = 14.59011781029012
$ws.Range("F12").Value = 27.04640138421247
$ws.Range("G12").Value = 24.55544665649753
$ws.Range("H12").Value = 13.24050657589022
$ws.Range("I12").Value = 18.44183201445918
$ws.Range("J12").Value = 9.719662556557159
$ws.Range("K12").Value = 19.80897884561086
$ws.Range("O12").Value = 19.59005817820074

# Row 13
$ws.Range("C13").Value = 4.600017083388022
$ws.Range("D13").Value = 10.77685577334964
$ws.Range("E13").Value = 14.58386464269635
$ws.Range("F13").Value = 27.05086369458545
$ws.Range("G13").Value = 24.56413511005557
$ws.Range("H13").Value = 13.24432704065868
$ws.Range("I13").Value = 18.44669073354749
$ws.Range("J13").Value = 9.719455600413051
$ws.Range("K13").Value = 19.76753650873419
$ws.Range("O13").Value = 19.59678236113111

# Row 14
$ws.Range("C14").Value = 4.572764287055607
$ws.Range("D14").Value = 10.75705670236357
$ws.Range("E14").Value = 14.5635198208239
$ws.Range("F14").Value = 27.0656538513017
$ws.Range("G14").Value = 24.59280557913457
$ws.Range("H14").Value = 13.25685907811916
$ws.Range("I14").Value = 18.46268797665179
$ws.Range("J14").Value = 9.718821116583783
$ws.Range("K14").Value = 19.63173303980304
$ws.Range("O14").Value = 19.61886632377392

# Row 15
$ws.Range("C15").Value = 4.555977299811577
$ws.Range("D15").Value = 10.7449360278351
$ws.Range("E15").Value = 14.55109537194022
$ws.Range("F15").Value = 27.07489564541219
$ws.Range("G15").Value = 24.61062439411488
$ws.Range("H15").Value = 13.2645911467296
$ws.Range("I15").Value = 18.47260336423964
$ws.Range("J15").Value = 9.718463572523344
$ws.Range("K15").Value = 19.54805050748926
$ws.Range("O15").Value = 19.63251235401969

# Row 16
$ws.Range("C16").Value = 4.458351105210199
$ws.Range("D16").Value = 10.67559698302267
$ws.Range("E16").Value = 14.480483342722
$ws.Range("F16").Value = 27.1306449649577
$ws.Range("G16").Value = 24.71666002103226
$ws.Range("H16").Value = 13.30975595413149
$ws.Range("I16").Value = 18.53120721125502
$ws.Range("J16").Value = 9.716891912467448
$ws.Range("K16").Value = 19.06090615346012
$ws.Range("O16").Value = 19.71253181396792

# Row 17
$ws.Range("C17").Value = 4.397201164131139
$ws.Range("D17").Value = 10.63318279616492
$ws.Range("E17").Value = 14.43770565461331
$ws.Range("F17").Value = 27.16732536853569
$ws.Range("G17").Value = 24.78518744576975
$ws.Range("H17").Value = 13.33822608562426
$ws.Range("I17").Value = 18.56874455653691
$ws.Range("J17").Value = 9.716354327038299
$ws.Range("K17").Value = 18.75534643345901
$ws.Range("O17").Value = 19.76324067310474

# Row 18
$ws.Range("C18").Value = 4.361572250586491
$ws.Range("D18").Value = 10.60883897869702
$ws.Range("E18").Value = 14.41330609547759
$ws.Range("F18").Value = 27.18933109837121
$ws.Range("G18").Value = 24.82587173498049
$ws.Range("H18").Value = 13.35488161497848
$ws.Range("I18").Value = 18.59091582624804
$ws.Range("J18").Value = 9.716202146449717
$ws.Range("K18").Value = 18.57715676884601
$ws.Range("O18").Value = 19.79300061397988

# Row 19
$ws.Range("C19").Value = 4.349430985166711
$ws.Range("D19").Value = 10.60060607409577
$ws.Range("E19").Value = 14.40508064023398
$ws.Range("F19").Value = 27.19693751851065
$ws.Range("G19").Value = 24.83986385326178
$ws.Range("H19").Value = 13.36056901451366
$ws.Range("I19").Value = 18.59852219843835
$ws.Range("J19").Value = 9.716177612904934
$ws.Range("K19").Value = 18.51640851515998
$ws.Range("O19").Value = 19.80317860897219

# Row 20
$ws.Range("C20").Value = 4.403758121467773
$ws.Range("D20").Value = 10.63769265388599
$ws.Range("E20").Value = 14.4422383354426
$ws.Range("F20").Value = 27.16332661426746
$ws.Range("G20").Value = 24.77776103084249
$ws.Range("H20").Value = 13.33516638249757
$ws.Range("I20").Value = 18.5646884884872
$ws.Range("J20").Value = 9.716395308079766
$ws.Range("K20").Value = 18.78812693687823
$ws.Range("O20").Value = 19.75778116244078

# Row 21
$ws.Range("C21").Value = 4.580788682363524
$ws.Range("D21").Value = 10.76287078230329
$ws.Range("E21").Value = 14.56948784248892
$ws.Range("F21").Value = 27.06127159528333
$ws.Range("G21").Value = 24.58433069116108
$ws.Range("H21").Value = 13.25316647840282
$ws.Range("I21").Value = 18.45796488522577
$ws.Range("J21").Value = 9.719001001218599
$ws.Range("K21").Value = 19.67172592402779
$ws.Range("O21").Value = 19.61235492707079

# Row 22
$ws.Range("C22").Value = 4.692302209889373
$ws.Range("D22").Value = 10.84500495869327
$ws.Range("E22").Value = 14.65434143029663
$ws.Range("F22").Value = 27.00272767269061
$ws.Range("G22").Value = 24.46940574810386
$ws.Range("H22").Value = 13.20207342610158
$ws.Range("I22").Value = 18.39342886886594
$ws.Range("J22").Value = 9.722096459074265
$ws.Range("K22").Value = 20.22694474741921
$ws.Range("O22").Value = 19.52263070203141

# Row 23
$ws.Range("C23").Value = 4.633173482399677
$ws.Range("D23").Value = 10.80114547752781
$ws.Range("E23").Value = 14.60890568565643
$ws.Range("F23").Value = 27.03322428692659
$ws.Range("G23").Value = 24.52968249750643
$ws.Range("H23").Value = 13.22911397405208
$ws.Range("I23").Value = 18.42739382803241
$ws.Range("J23").Value = 9.720317247627126
$ws.Range("K23").Value = 19.93267441547862
$ws.Range("O23").Value = 19.57002969490355

# Row 24
$ws.Range("C24").Value = 4.400795192999888
$ws.Range("D24").Value = 10.63565361902672
$ws.Range("E24").Value = 14.44018850495133
$ws.Range("F24").Value = 27.16513159331646
$ws.Range("G24").Value = 24.78111450633261
$ws.Range("H24").Value = 13.33654877821261
$ws.Range("I24").Value = 18.56652039777637
$ws.Range("J24").Value = 9.716376291810789
$ws.Range("K24").Value = 18.77331471204736
$ws.Range("O24").Value = 19.76024751759618

# Row 25
$ws.Range("C25").Value = 4.132758165425978
$ws.Range("D25").Value = 10.45920127665603
$ws.Range("E25").Value = 14.2661139027255
$ws.Range("F25").Value = 27.34208971294769
$ws.Range("G25").Value = 25.10073183432582
$ws.Range("H25").Value = 13.46314859058611
$ws.Range("I25").Value = 18.7387704582477
$ws.Range("J25").Value = 9.718130776756897
$ws.Range("K25").Value = 17.42995494914596
$ws.Range("O25").Value = 19.98810188452471
